$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new sheets: META (before CAPABILITIES) and DESCRIPTIONS (after).
#    Always re-fetch sheets by name -- the COM object returned by Item() in
#    this runtime tracks the *position*, not a stable identity.
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($wb.Worksheets.Item("CAPABILITIES"))
$metaSheet.Name = "META"

$descSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item("CAPABILITIES"))
$descSheet.Name = "DESCRIPTIONS"

# ---------------------------------------------------------------------------
# 2. CAPABILITIES: the "current_unit" column (P) changes from "countable" to
#    "percent" for every data row. A few rows (7, 8, 19-22) did not carry the
#    grey banding style that the rest of column P has; pick it up by copying
#    formats from a cell that already has it before overwriting the values.
# ---------------------------------------------------------------------------
$cap = $wb.Worksheets.Item("CAPABILITIES")

$cap.Range("P2").Copy()
$cap.Range("P7").PasteSpecial(-4122)
$cap.Range("P8").PasteSpecial(-4122)
$cap.Range("P19").PasteSpecial(-4122)
$cap.Range("P20").PasteSpecial(-4122)
$cap.Range("P21").PasteSpecial(-4122)
$cap.Range("P22").PasteSpecial(-4122)

for ($r = 2; $r -le 23; $r++) {
    $cap.Range("P$r").Value = "percent"
}

# New custom width picked up on column P (16).
$cap.Columns.Item(16).ColumnWidth = 9.33

# Selection moves off A25 (no longer a selected/active sheet) to A3.
$cap.Range("A3").Select()

# ---------------------------------------------------------------------------
# 3. DESCRIPTIONS: two-language free text sheet, keyed by capability_id.
# ---------------------------------------------------------------------------
$desc = $wb.Worksheets.Item("DESCRIPTIONS")

$desc.Columns.Item(1).ColumnWidth = 35.0
$desc.Columns.Item(2).ColumnWidth = 105.0
$desc.Columns.Item(3).ColumnWidth = 113.33333333333333

$desc.Range("A1").Value = "capability_id"
$desc.Range("B1").Value = "en"
$desc.Range("C1").Value = "fr"

$desc.Range("A2").Value = "DF Regt"
$desc.Range("B2").Value = "Lorem ipsum dolor sit amet, consectetur adipiscing elit. Proin ac est et ante sagittis volutpat vel non nunc. Nunc gravida est ex, elementum tincidunt enim accumsan at. Fusce vel tempor lorem. Nunc a nulla nulla. Sed faucibus, erat et iaculis dapibus, ligula lacus iaculis tellus, a tempus ex urna vel magna. Donec eu metus eu ante dignissim blandit non vitae mi. In metus massa, congue sed massa nec, aliquam gravida nisi. Nam ornare, ante ac tempor maximus, velit enim rhoncus ligula, vitae pretium neque nunc eget erat. Nullam aliquet lectus vitae massa eleifend, vitae mattis mi sollicitudin. Maecenas sed mollis sem. Proin et magna ultrices, ultricies mauris at, elementum risus. Nullam aliquam fermentum mi tincidunt pulvinar."
$desc.Range("C2").Value = "Praesent eu sapien maximus, vulputate sapien vel, iaculis mauris. Suspendisse erat odio, tristique vel ligula sed, rutrum egestas mauris. Donec non maximus sem. Quisque ornare, ligula ut bibendum vestibulum, nulla dolor volutpat quam, congue finibus dui tellus non metus. Mauris euismod, metus vitae feugiat bibendum, nisi leo condimentum mi, et ultricies erat felis vitae purus. Donec tortor nulla, laoreet in tristique non, pharetra ultrices tortor. Nulla tempus tortor non sapien tincidunt, eu imperdiet lacus dignissim. Proin id orci libero. Phasellus mattis tempor velit, non vestibulum libero fermentum sed. Curabitur tincidunt ipsum non luctus dapibus. Curabitur et vehicula justo. Mauris ut purus pharetra, porttitor ex vel, venenatis urna. Vivamus ornare urna eu arcu euismod, vel volutpat mauris pulvinar. Nam elementum dui at velit sollicitudin hendrerit."

$desc.Range("A3").Value = "IDF Regts"
$desc.Range("B3").Value = "nterdum et malesuada fames ac ante ipsum primis in faucibus. Nam quis gravida lacus, quis gravida neque. Morbi hendrerit turpis ac mollis venenatis. Donec ac ullamcorper nisl. In hac habitasse platea dictumst. Cras tempor posuere bibendum. Vivamus consectetur accumsan dolor quis eleifend. Quisque quis magna mauris."

$desc.Range("B1:C2").WrapText = $true
$desc.Range("B3").WrapText = $true

$desc.Rows.Item(2).RowHeight = 120
$desc.Rows.Item(3).RowHeight = 60

$desc.Range("C3").Select()

# ---------------------------------------------------------------------------
# 4. Make DESCRIPTIONS the active/visible tab (sets workbookView activeTab
#    and moves tabSelected onto this sheet's sheetView).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("DESCRIPTIONS").Activate()
